$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 55 - 168. Excel Sheet Column Title (Easy, finished)
# Template: row 48 ("326. Power of Three") has the identical style pattern
# (A=32 B=4 C=22 D=6 E=23 F=23 G=44) that row 55 needs.
# ---------------------------------------------------------------------------
$ws.Range("A48:G48").Copy($ws.Range("A55:G55"))
$ws.Cells.Item(55, 1).Value = "168. Excel Sheet Column Title"
$ws.Cells.Item(55, 2).Value = "Easy"
$ws.Cells.Item(55, 3).Value = "https://leetcode.com/problems/excel-sheet-column-title/"
$ws.Cells.Item(55, 4).Value = 44558
$ws.Cells.Item(55, 5).Value = "数学，进制转化"
$ws.Cells.Item(55, 6).Value = "注意进制起始位的细节"
$ws.Cells.Item(55, 7).Value = "未复习"

# ---------------------------------------------------------------------------
# Row 56 - 67. Add Binary (Easy, finished)
# Same template as row 55, except column F keeps the "English-first" style
# (style of column B, Times New Roman) because the note begins with the
# Latin word "StringBuffer" before switching to Chinese.
# ---------------------------------------------------------------------------
$ws.Range("A48:G48").Copy($ws.Range("A56:G56"))
$ws.Range("B48").Copy($ws.Range("F56"))
$ws.Cells.Item(56, 1).Value = "67. Add Binary"
$ws.Cells.Item(56, 2).Value = "Easy"
$ws.Cells.Item(56, 3).Value = "https://leetcode.com/problems/add-binary/"
$ws.Cells.Item(56, 4).Value = 44558
$ws.Cells.Item(56, 5).Value = "数学，字符串加法"
$ws.Cells.Item(56, 7).Value = "未复习"

$f56 = $ws.Cells.Item(56, 6)
$f56.Value = "StringBuffer翻转函数；长度不同时按位相加的写法"
$prefixLen = "StringBuffer".Length
$totalLen = $f56.Value.Length
$chineseRun = $f56.Characters($prefixLen + 1, $totalLen - $prefixLen)
$chineseRun.Font.Name = "宋体"
$chineseRun.Font.Size = 11

# ---------------------------------------------------------------------------
# Row 57 - 540. Single Element in a Sorted Array (Medium, reviewed)
# Template: row 44 ("451. Sort Characters By Frequency") matches the style
# pattern row 57 needs (A=34 B=18 C=19 D=20 E=21 F=31 G=28 H=42); only
# column F uses the wrap-capable variant of that style, taken from row 35.
# ---------------------------------------------------------------------------
$ws.Range("A44:H44").Copy($ws.Range("A57:H57"))
$ws.Range("F35").Copy($ws.Range("F57"))
$ws.Cells.Item(57, 1).Value = "540. Single Element in a Sorted Array"
$ws.Cells.Item(57, 2).Value = "Medium"
$ws.Cells.Item(57, 3).Value = "https://leetcode.com/problems/single-element-in-a-sorted-array/"
$ws.Cells.Item(57, 4).Value = 44470
$ws.Cells.Item(57, 5).Value = "二分法"
$ws.Cells.Item(57, 6).Value = "奇数长度数组，mid左右元素个数相同；通过mid和mid+1元素比较锁定单个元素位置在左边还是右边"
$ws.Cells.Item(57, 7).Value = 44558
$ws.Cells.Item(57, 8).Value = "√"

# ---------------------------------------------------------------------------
# Hyperlinks for the 3 new "link" cells. Hyperlinks.Add re-applies its own
# formatting and clobbers the cell style with a freshly duplicated one, so
# re-paste the already-correct format from the template cell right after.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C55"), "https://leetcode.com/problems/excel-sheet-column-title/") | Out-Null
$ws.Range("C48").Copy()
$ws.Range("C55").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C56"), "https://leetcode.com/problems/add-binary/") | Out-Null
$ws.Range("C48").Copy()
$ws.Range("C56").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C57"), "https://leetcode.com/problems/single-element-in-a-sorted-array/") | Out-Null
$ws.Range("C44").Copy()
$ws.Range("C57").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet view bookkeeping: scroll so the new rows are visible and move the
# active selection the way the author left it.
# ---------------------------------------------------------------------------
$ws.Range("A52").Select()
$ws.Range("I61").Select()
